$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking price text cells so they stay as text (match source formatting)
$textCells = @("D4", "D5", "D6", "D8", "D9", "D11", "D12", "D16", "D17", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D45", "D46", "D47", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell updates row by row
# Row 2
$ws.Range("D2").Value = "72.169.37"
$ws.Range("E2").Value = "  -0.46%  "

# Row 3
$ws.Range("D3").Value = "3.907.38"
$ws.Range("E3").Value = "  -1.77%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "602.00"
$ws.Range("E5").Value = "  +2.33%  "

# Row 6
$ws.Range("D6").Value = "168.58"

# Row 7
$ws.Range("E7").Value = "  -0.56%  "

# Row 8
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").Value = "0.769"
$ws.Range("E9").Value = "  +3.18%  "

# Row 10
$ws.Range("E10").Value = "  +8.26%  "

# Row 11
$ws.Range("D11").Value = "54.55"
$ws.Range("E11").Value = "  +2.89%  "

# Row 12
$ws.Range("D12").Value = "0.0000325"
$ws.Range("E12").Value = "  +2.67%  "

# Row 13
$ws.Range("E13").Value = "  +6.27%  "

# Row 14
$ws.Range("D14").Value = "4.537.23"
$ws.Range("E14").Value = "  -1.77%  "

# Row 15
$ws.Range("D15").Value = "3.921.27"
$ws.Range("E15").Value = "  -1.71%  "

# Row 16
$ws.Range("D16").Value = "21.17"
$ws.Range("E16").Value = "  +3.83%  "

# Row 17
$ws.Range("D17").Value = "14.01"
$ws.Range("E17").Value = "  +0.38%  "

# Row 18
$ws.Range("E18").Value = "  -4.90%  "

# Row 19
$ws.Range("D19").Value = "72.063.15"
$ws.Range("E19").Value = "  -0.61%  "

# Row 20
$ws.Range("E20").Value = "  -1.67%  "

# Row 21
$ws.Range("D21").Value = "438.29"
$ws.Range("E21").Value = "  +2.32%  "

# Row 22
$ws.Range("D22").Value = "4.79"
$ws.Range("E22").Value = "  +1.70%  "

# Row 23
$ws.Range("D23").Value = "94.69"
$ws.Range("E23").Value = "  -0.91%  "

# Row 24
$ws.Range("E24").Value = "  -3.41%  "

# Row 25
$ws.Range("D25").Value = "13.93"
$ws.Range("E25").Value = "  -1.86%  "

# Row 26
$ws.Range("D26").Value = "4.17"
$ws.Range("E26").Value = "  -7.28%  "

# Row 27
$ws.Range("D27").Value = "11.06"
$ws.Range("E27").Value = "  -1.53%  "

# Row 28
$ws.Range("D28").Value = "5.95"
$ws.Range("E28").Value = "  +0.39%  "

# Row 29
$ws.Range("D29").Value = "10.26"
$ws.Range("E29").Value = "  -2.52%  "

# Row 30
$ws.Range("D30").Value = "35.35"
$ws.Range("E30").Value = "  -2.37%  "

# Row 31
$ws.Range("D31").Value = "7.98"
$ws.Range("E31").Value = "  +2.31%  "

# Row 32
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "52.18"
$ws.Range("E32").Value = "  +4.53%  "

# Row 33
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "13.71"
$ws.Range("E33").Value = "  +2.12%  "

# Row 34
$ws.Range("E34").Value = "  -3.25%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0987"
$ws.Range("E35").Value = "  +16.19%  "

# Row 36
$ws.Range("D36").Value = "68.44"
$ws.Range("E36").Value = "  -0.09%  "

# Row 37
$ws.Range("D37").Value = "619.43"
$ws.Range("E37").Value = "  -9.00%  "

# Row 38
$ws.Range("D38").Value = "0.423"
$ws.Range("E38").Value = "  -3.07%  "

# Row 39
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.10%  "

# Row 40
$ws.Range("D40").Value = "3.33"
$ws.Range("E40").Value = "  +0.98%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.08%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.143"
$ws.Range("E42").Value = "  -1.89%  "

# Row 43
$ws.Range("E43").Value = "  +41.68%  "

# Row 44
$ws.Range("E44").Value = "  -2.43%  "

# Row 45
$ws.Range("D45").Value = "10.34"
$ws.Range("E45").Value = "  -6.72%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.145"
$ws.Range("E46").Value = "  -1.96%  "

# Row 47
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "2.65"
$ws.Range("E47").Value = "  -3.71%  "

# Row 48
$ws.Range("E48").Value = "  -15.57%  "

# Row 49
$ws.Range("D49").Value = "3.33"
$ws.Range("E49").Value = "  -0.69%  "

# Row 50
$ws.Range("D50").Value = "2.888.06"
$ws.Range("E50").Value = "  +4.09%  "

# Row 51
$ws.Range("E51").Value = "  +3.48%  "
